# Correct Problem Statement and Success Criteria
# Applies targeted text edits to the single slide of the Nordic Sensing Co
# "Problem Identification" deck:
#   1. Success-criteria bullet: clarify what "failures" means (15% failure rate).
#   2. Problem-statement question: turn the "how" question into a "should we"
#      statement, and spell out the two concrete remediation options
#      (shutting down a factory / stopping a supplier).
#   3. Scope blurb: scope is the supply chain, not the manufacturing process.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "Success will be measured in our capacity to: ... Detect the origin of
#    the failures ..." -> "... Detect the origin of the 15% failure rate"
# ---------------------------------------------------------------------------
$shpSuccess = $s.Shapes.Item("Google Shape;35;p1")
$trSuccess = $shpSuccess.TextFrame.TextRange

$oldText1 = "Detect the origin of the failures"
$newText1 = "Detect the origin of the 15% failure rate"
$idx1 = $trSuccess.Text.IndexOf($oldText1)
if ($idx1 -ge 0) {
    $trSuccess.Characters($idx1 + 1, $oldText1.Length).Text = $newText1
}

# ---------------------------------------------------------------------------
# 2) Problem statement headline shape: rewrite the question and extend it
#    with the two remediation options.
# ---------------------------------------------------------------------------
$shpQuestion = $s.Shapes.Item("Google Shape;48;p1")
$trQuestion = $shpQuestion.TextFrame.TextRange

# 2a) "How can we determine the origin ... NSC's " -> "Determine if the origin ... NSC's "
#      (PowerPoint's TextRange.Text always reports a plain apostrophe, but the
#      run actually stores a right single quotation mark (U+2019); write the
#      real character back so the underlying XML keeps matching punctuation.)
$rsquo = [char]0x2019
$oldLead = "How can we determine the origin of the 15% failure rate in the manufacturing process of the NSC's "
$newLead = "Determine if the origin of the 15% failure rate in the manufacturing process of the NSC" + $rsquo + "s "
$idx2 = $trQuestion.Text.IndexOf($oldLead)
if ($idx2 -ge 0) {
    $trQuestion.Characters($idx2 + 1, $oldLead.Length).Text = $newLead
}

# 2b) " sensor by April the 3" -> " sensor should be solved by shutting down a
#      factory or stop buying certain parts from a supplier by April the 3"
#      (the "rd" superscript + trailing "? " runs that already follow are left
#      untouched).
$oldTail = " sensor by April the 3"
$newTail = " sensor should be solved by shutting down a factory or stop buying certain parts from a supplier by April the 3"
$idx3 = $trQuestion.Text.IndexOf($oldTail)
if ($idx3 -ge 0) {
    $tailRange = $trQuestion.Characters($idx3 + 1, $oldTail.Length)
    $tailRange.Text = $newTail

    # The replaced span now reads (all one run, inheriting $oldTail's format):
    #   " sensor should be solved by shutting down a factory or stop buying
    #    certain parts from a supplier by April the 3"
    # Split out the middle clause ("g down a factory ... from a supplier ")
    # into its own lightly-formatted runs, matching how it was typed fresh
    # (no explicit Arial/size/color overrides), while the leading
    # (" sensor should be solved by shuttin") and trailing ("by April the 3")
    # pieces keep the original Arial/14pt/black formatting.
    $runStart = $tailRange.Start

    $keepLead = " sensor should be solved by shuttin"
    $midPart1 = "g down a factory or stop buying certain parts "
    $midPart2 = "from a supplier "
    $keepTrail = "by April the 3"

    $midStart = $runStart + $keepLead.Length

    # Touching .Font on each sub-range forces PowerPoint to break it out into
    # its own <a:r> run (matching the 4-way run split in the target), even
    # though the effective Arial/14pt/black formatting ends up the same as
    # the surrounding text.
    $part1Range = $trQuestion.Characters($midStart, $midPart1.Length)
    $part1Range.Font.Size = $part1Range.Font.Size

    $part2Range = $trQuestion.Characters($midStart + $midPart1.Length, $midPart2.Length)
    $part2Range.Font.Size = $part2Range.Font.Size

    $trailStart = $midStart + $midPart1.Length + $midPart2.Length
    $trailRange = $trQuestion.Characters($trailStart, $keepTrail.Length)
    $trailRange.Font.Size = $trailRange.Font.Size
}

# ---------------------------------------------------------------------------
# 3) Scope blurb: "Our scope will be on the manufacturing process for the " ->
#    "Our scope will be on the supply chain for the "
# ---------------------------------------------------------------------------
$shpScope = $s.Shapes.Item("Google Shape;36;p1")
$trScope = $shpScope.TextFrame.TextRange

$oldScope = "Our scope will be on the manufacturing process for the "
$newScope = "Our scope will be on the supply chain for the "
$idx4 = $trScope.Text.IndexOf($oldScope)
if ($idx4 -ge 0) {
    $trScope.Characters($idx4 + 1, $oldScope.Length).Text = $newScope
}
